# This script applies the stock-report corrections described in the
# commit diff: quantity (F) and value (G) corrections for several stock
# lines, a couple of rows whose Code/Rate/Qty/Value (B/E/F/G) were
# swapped with their neighboring row, and the resulting Sub Total /
# Grand Total rollups (column B) that depend on those line values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70
$ws.Range("F70").Value = 50
$ws.Range("G70").Value = 6747.5

# Row 78
$ws.Range("F78").Value = 53
$ws.Range("G78").Value = 3015.7

# Row 90
$ws.Range("B90").Value = 217382.27

# Row 196
$ws.Range("F196").Value = 4
$ws.Range("G196").Value = 458.24

# Row 205
$ws.Range("F205").Value = 32
$ws.Range("G205").Value = 12068.48

# Row 215
$ws.Range("F215").Value = 4
$ws.Range("G215").Value = 1692.84

# Row 216
$ws.Range("B216").Value = 58123.36

# Row 227
$ws.Range("B227").Value = 55373
$ws.Range("E227").Value = 163.62
$ws.Range("F227").Value = -94
$ws.Range("G227").Value = -13562.32

# Row 228
$ws.Range("B228").Value = 63520
$ws.Range("E228").Value = 153.4
$ws.Range("F228").Value = 68
$ws.Range("G228").Value = 9811.04

# Row 234
$ws.Range("F234").Value = 48
$ws.Range("G234").Value = 2463.36

# Row 260
$ws.Range("B260").Value = 227827.37

# Row 364
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65

# Row 365
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45

# Row 366
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9

# Row 367
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29

# Row 372
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46

# Row 373
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91

# Row 375
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68

# Row 376
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32

# Row 380
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65

# Row 381
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945

# Row 382
$ws.Range("B382").Value = 45702
$ws.Range("E382").Value = 31.43
$ws.Range("F382").Value = -215
$ws.Range("G382").Value = -5654.5

# Row 383
$ws.Range("B383").Value = 64919
$ws.Range("E383").Value = 27.97
$ws.Range("F383").Value = 61
$ws.Range("G383").Value = 1604.3

# Row 385
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98

# Row 386
$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55

# Row 473
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79

# Row 474
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 109
$ws.Range("G474").Value = 3578.47

# Row 553
$ws.Range("F553").Value = 7
$ws.Range("G553").Value = 529.76

# Row 560
$ws.Range("B560").Value = 35162.54

# Row 619
$ws.Range("B619").Value = 2225667.9

# Row 620
$ws.Range("B620").Value = 2225667.9
